$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 66, shifting the existing
# rows 66-67 down to 68-69 (their contents remain unchanged).
$ws.Range("A66:A67").EntireRow.Insert()

# New row 66 - weekly update for "Primera" quality
$ws.Range("A66").Value = 7
$ws.Range("B66").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C66").Value = "Ñuble"
$ws.Range("D66").Value = 45239
$ws.Range("E66").Value = 16
$ws.Range("F66").Value = 300000000
$ws.Range("G66").Value = "Espárragos"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 600
$ws.Range("K66").Value = 1200
$ws.Range("L66").Value = 1300
$ws.Range("M66").Value = 1250
$ws.Range("N66").Value = "$/kilo"
$ws.Range("O66").Value = "Región de Ñuble"
$ws.Range("P66").Value = 1250
$ws.Range("Q66").Value = 1
$ws.Range("R66").Value = "Hortaliza"

# New row 67 - weekly update for "Segunda" quality
$ws.Range("A67").Value = 7
$ws.Range("B67").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C67").Value = "Ñuble"
$ws.Range("D67").Value = 45239
$ws.Range("E67").Value = 16
$ws.Range("F67").Value = 300000000
$ws.Range("G67").Value = "Espárragos"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Segunda"
$ws.Range("J67").Value = 600
$ws.Range("K67").Value = 1000
$ws.Range("L67").Value = 1000
$ws.Range("M67").Value = 1000
$ws.Range("N67").Value = "$/kilo"
$ws.Range("O67").Value = "Región de Ñuble"
$ws.Range("P67").Value = 1000
$ws.Range("Q67").Value = 1
$ws.Range("R67").Value = "Hortaliza"
